# fix bug in import/export
# Move the two values that used to live in B6:B7 over to F3:F4,
# and fix up their text ("3rqwe" / "eqwrewq" instead of "test" / 1231).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old data out of B6:B7 entirely.
$ws.Range("B6:B7").ClearContents()

# Write the corrected values into F3:F4.
$ws.Range("F3").Value = "3rqwe"
$ws.Range("F4").Value = "eqwrewq"

# Leave the selection on the last-edited cell, like the source file.
$ws.Range("F4").Select()
